# Revert capacity charts to show kilowatts on the y-axis.
# The underlying sheet data (and the chart's cached values) were stored in
# watts; convert the affected cells to kilowatts (divide by 1000) and
# update the axis title / number formats to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Worksheet data: convert Solar (column E) and the one Energy Storage
# (column C) outlier from watts to kilowatts. ---
$ws.Range("E16").Value = 4
$ws.Range("E19").Value = 13.6
$ws.Range("E20").Value = 6.7
$ws.Range("E21").Value = 21.5
$ws.Range("E22").Value = 20.1
$ws.Range("E23").Value = 100.6
$ws.Range("E24").Value = 96.65000000000001
$ws.Range("E25").Value = 197.341
$ws.Range("E26").Value = 130.58
$ws.Range("C25").Value = 3.85

# --- Number format for the data grid now needs one decimal place, since
# values are fractional kilowatts. ---
$ws.Range("B2:G26").NumberFormat = "#,##0.0"

# --- Chart: update axis title and number format on the value axis. ---
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$valAx = $chart.Axes(2)
$valAx.AxisTitle.Text = "Kilowatts (kW)"
$valAx.TickLabels.NumberFormat = "#,##0"
